# Auto-generated edit script: applies the 2026-01-19 20:12:40 scrape refresh
# to all three schedule sheets (row 2/3 header text + changed/appended data rows).
$wb = $excel.ActiveWorkbook

# ---- Sheet 1: LP1912 ----
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 1).Value = "Última actualización: 20:12:40"
$ws.Cells.Item(3, 1).Value = "Total filas: 480"

$rows = @(
    @(79, "06:25:43", "08:23", "215B_EL PATO", 118, "LP1912"),
    @(80, "06:55:02", "08:23", "16_P MOR-SANTA ANA", 88, "LP1912"),
    @(110, "08:56:29", "09:34", "23_HERNANDEZ", 38, "LP1912"),
    @(111, "08:56:29", "09:34", "16_SANTA ANA", 38, "LP1912"),
    @(215, "11:30:45", "13:07", "10_OLMOS", 97, "LP1912"),
    @(216, "12:21:08", "13:07", "16_P MOR-SANTA ANA", 46, "LP1912"),
    @(335, "16:13:37", "17:04", "11_ETCHEVERRY", 51, "LP1912"),
    @(336, "16:13:37", "17:04", "23_HERNANDEZ", 51, "LP1912"),
    @(337, "15:22:17", "17:04", "215A_EL PATO", 102, "LP1912"),
    @(364, "15:53:28", "17:52", "81_EL PELIGRO", 119, "LP1912"),
    @(366, "17:38:19", "17:52", "11_ETCHEVERRY", 14, "LP1912"),
    @(410, "18:33:37", "19:16", "15_ABASTO", 43, "LP1912"),
    @(411, "17:51:34", "19:16", "27_EL RETIRO", 85, "LP1912"),
    @(448, "20:12:40", "20:12", "10_OLMOS", 0, "LP1912"),
    @(449, "19:39:04", "20:13", "23_HERNANDEZ", 34, "LP1912"),
    @(450, "20:12:40", "20:13", "16_SANTA ANA", 1, "LP1912"),
    @(451, "18:33:37", "20:20", "26_HERNANDEZ", 107, "LP1912"),
    @(452, "18:47:08", "20:21", "26_HERNANDEZ", 94, "LP1912"),
    @(453, "18:47:08", "20:21", "16_P MOR-SANTA ANA", 94, "LP1912"),
    @(454, "19:52:18", "20:22", "16_SANTA ANA", 30, "LP1912"),
    @(455, "18:33:37", "20:22", "11_ETCHEVERRY", 109, "LP1912"),
    @(456, "18:33:37", "20:23", "215A_EL PATO", 110, "LP1912"),
    @(457, "18:47:08", "20:24", "14_ABASTO", 97, "LP1912"),
    @(458, "18:47:08", "20:24", "215A_EL PATO", 97, "LP1912"),
    @(459, "18:53:41", "20:28", "16_P MOR-SANTA ANA", 95, "LP1912"),
    @(460, "18:53:41", "20:30", "14_ABASTO", 97, "LP1912"),
    @(461, "18:33:37", "20:31", "225_GOMEZ", 118, "LP1912"),
    @(462, "19:52:18", "20:32", "16_P MOR-SANTA ANA", 40, "LP1912"),
    @(463, "18:47:08", "20:39", "11_ETCHEVERRY", 112, "LP1912"),
    @(464, "19:13:07", "20:40", "16_P MOR-SANTA ANA", 87, "LP1912"),
    @(465, "19:13:07", "20:44", "11_ETCHEVERRY", 91, "LP1912"),
    @(466, "20:12:40", "20:47", "16_P MOR-SANTA ANA", 35, "LP1912"),
    @(467, "19:13:07", "20:52", "15_ABASTO", 99, "LP1912"),
    @(468, "19:52:18", "20:52", "23_HERNANDEZ", 60, "LP1912"),
    @(469, "20:12:40", "20:56", "27_EL RETIRO", 44, "LP1912"),
    @(470, "19:52:18", "20:56", "10_OLMOS", 64, "LP1912"),
    @(471, "19:13:07", "20:57", "27_EL RETIRO", 104, "LP1912"),
    @(472, "19:13:07", "21:04", "84_COLONIA URQUIZA-ESC 49", 111, "LP1912"),
    @(473, "20:12:40", "21:07", "215B_EL PATO", 55, "LP1912"),
    @(474, "19:13:07", "21:08", "215B_EL PATO", 115, "LP1912"),
    @(475, "20:12:40", "21:08", "16_P MOR-167 Y 521", 56, "LP1912"),
    @(476, "19:39:04", "21:16", "14_ABASTO", 97, "LP1912"),
    @(477, "20:12:40", "21:20", "26_HERNANDEZ", 68, "LP1912"),
    @(478, "19:39:04", "21:21", "26_HERNANDEZ", 102, "LP1912"),
    @(479, "19:39:04", "21:23", "10_OLMOS", 104, "LP1912"),
    @(480, "20:12:40", "21:37", "17_ROMERO", 85, "LP1912"),
    @(481, "19:52:18", "21:38", "14_ABASTO", 106, "LP1912"),
    @(482, "19:39:04", "21:38", "17_ROMERO", 119, "LP1912"),
    @(483, "19:52:18", "21:47", "215A_EL PATO", 115, "LP1912"),
    @(484, "20:12:40", "21:52", "10_OLMOS", 100, "LP1912"),
    @(485, "20:12:40", "22:07", "17_ROMERO", 115, "LP1912")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}

# ---- Sheet 2: LP1912-215 ----
$ws = $wb.Worksheets.Item(2)

$ws.Cells.Item(2, 1).Value = "Última actualización: 20:12:40"
$ws.Cells.Item(3, 1).Value = "Total filas: 49"

$rows = @(
    @(52, "20:12:40", "21:07", "215B_EL PATO", 55, "LP1912"),
    @(53, "19:13:07", "21:08", "215B_EL PATO", 115, "LP1912"),
    @(54, "19:52:18", "21:47", "215A_EL PATO", 115, "LP1912")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}

# ---- Sheet 3: 6203-6173 ----
$ws = $wb.Worksheets.Item(3)

$ws.Cells.Item(2, 1).Value = "Última actualización: 20:12:40"
$ws.Cells.Item(3, 1).Value = "Total filas: 62"

$rows = @(
    @(65, "20:12:40", "21:28", "215C_LA PLATA", 76, "L6203"),
    @(66, "19:39:04", "21:29", "215C_LA PLATA", 110, "L6203"),
    @(67, "20:12:40", "22:04", "215A_LA PLATA", 112, "L6173")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}

Write-Output "edit complete"